$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added at the top of the data table
# (row 27), pushing every existing data row (27-59) down by one
# (to 28-60). Insert a blank row at 27 first so the rest of the table
# shifts down intact, then populate the new row with the latest entry.
$ws.Rows("27:27").Insert()

$ws.Range("A27").Value = 1
$ws.Range("B27").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C27").Value = "Arica y Parinacota"
$ws.Range("D27").Value = 44645
$ws.Range("E27").Value = 15
$ws.Range("F27").Value = 100112027
$ws.Range("G27").Value = "Melón"
$ws.Range("H27").Value = "Calameño"
$ws.Range("I27").Value = "Primera"
$ws.Range("J27").Value = 80
$ws.Range("K27").Value = 13000
$ws.Range("L27").Value = 14000
$ws.Range("M27").Value = 13500
$ws.Range("N27").Value = "`$/caja 18 unidades"
$ws.Range("O27").Value = "Región de Arica y Parinacota"
$ws.Range("P27").Value = 750
$ws.Range("Q27").Value = 18
$ws.Range("R27").Value = "Hortaliza"
